# Översikt KARLSHAMN - automatic update of files.
#
# The sheet lists "Avverkningsanmälningar" (logging notifications). Each
# existing record's "Förändrad" (last changed) date in column C moves from
# 2023-09-06 (45175) to 2023-09-08 (45177), and one brand new record
# (A 41533-2023) is appended as row 329.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" date for every existing data row (2..328) from
#    45175 (2023-09-06) to 45177 (2023-09-08).
$ws.Range("C2:C328").Value = 45177

# 2) Row 328 picks up an explicit (default) row height in the new file.
$ws.Rows.Item(328).RowHeight = 15

# 3) Append the new record as row 329.
$ws.Range("A329").Value = "A 41533-2023"

$ws.Range("B329").Value = 45175
$ws.Range("B329").NumberFormat = "YYYY-MM-DD"

$ws.Range("C329").Value = 45177
$ws.Range("C329").NumberFormat = "YYYY-MM-DD"

$ws.Range("D329").Value = "BLEKINGE LÄN"
$ws.Range("E329").Value = "KARLSHAMN"

$ws.Range("G329").Value = 6.3
$ws.Range("H329").Value = 0
$ws.Range("I329").Value = 0
$ws.Range("J329").Value = 0
$ws.Range("K329").Value = 0
$ws.Range("L329").Value = 0
$ws.Range("M329").Value = 0
$ws.Range("N329").Value = 0
$ws.Range("O329").Value = 0
$ws.Range("P329").Value = 0
$ws.Range("Q329").Value = 0

# R is the wrapped "Artnamn" column; leave it blank but keep the same
# wrap-text style used throughout the rest of the column.
$ws.Range("R329").Value = ""
$ws.Range("R329").WrapText = $true
